# Re-colour the presentation's theme (ppt/theme/theme1.xml) from the
# "Integral" / "Red Violet" colour scheme to the default Office Theme
# colour scheme, by setting each of the 12 theme colours through the
# Design -> SlideMaster -> Theme -> ThemeColorScheme object model.
#
# NOTE: we deliberately never read/write the .Name property on the
# Design/Theme/ThemeColorScheme objects: those setters are no-ops in
# this host (they don't persist into the OOXML) and merely reading or
# writing them forces an unrelated, spurious re-serialisation of other
# slide parts. Sticking to the RGB setters keeps the edit limited to
# ppt/theme/theme1.xml, which is what actually carries the colour
# scheme data.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB = 0        # dk1      000000
$colors.Item(2).RGB = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB = 6968388  # dk2      44546A
$colors.Item(4).RGB = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB = 3243501  # accent2  ED7D31
$colors.Item(7).RGB = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB = 49407    # accent4  FFC000
$colors.Item(9).RGB = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456 # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink   0563C1
$colors.Item(12).RGB = 7491477  # folHlink 954F72
